$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting like "1.00" / "0.0000168"
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.887.36"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "3.230.23"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "575.26"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "175.52"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "3.228.07"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +4.93%  "
$ws.Range("D11").Value = "6.67"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "3.787.59"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "0.137"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "27.58"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "66.844.33"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").Value = "0.0000168"
$ws.Range("E17").Value = "  +4.06%  "
$ws.Range("D18").Value = "3.225.11"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "5.76"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "13.19"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "365.84"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").Value = "7.42"
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "69.79"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "0.0000119"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "0.503"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "3.349.78"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "1.95"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").Value = "5.55"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "22.35"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "6.71"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "168.03"
$ws.Range("E37").Value = "  +6.31%  "
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  +10.18%  "
$ws.Range("D41").Value = "26.56"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("D43").Value = "6.33"
$ws.Range("E43").Value = "  +5.06%  "
$ws.Range("D44").Value = "2.682.64"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "4.25"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").Value = "40.36"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("D47").Value = "0.0669"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "24.26"
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("D49").Value = "329.54"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("D50").Value = "0.0277"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("E51").Value = "  +0.93%  "
